# Adding saga audit plugin
#
# 1) Slide 14 ("Event driven architectures"): merge the two runs
#    "Where " + "sagas rule…" into a single run "Where sagas rule…<tab><tab>"
#    (keeping the first run's formatting).
# 2) Slide 3 ("Pre requisites" -> "Prerequisites"): retitle, and fill in
#    the previously-empty content placeholder with the prerequisites list.

$p = $ppt.ActivePresentation

# --- Edit 1: slide 14, Text Placeholder -------------------------------
$s14  = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)
$tr14 = $sh14.TextFrame.TextRange

# Keep the first run ("Where ") as-is and fold the rest of the text into
# it so the two original runs collapse into one run.
$run1 = $tr14.Characters(1, 6)
$rest = $tr14.Characters(7, $tr14.Length - 6)
$rest.Text = ""
[void]$run1.InsertAfter("sagas rule…`t`t")

# --- Edit 2: slide 3, Title ---------------------------------------------
$s3     = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "Prerequisites"

# --- Edit 3: slide 3, Content Placeholder -------------------------------
$content3 = $s3.Shapes.Item(2)
$tr3 = $content3.TextFrame.TextRange

# Build the last paragraph ("SQL Server" / " " / "+ Management tools")
# first, while it is the shape's only text, so each piece's language can
# be stamped correctly; then prepend the first two bullet paragraphs.
$tr3.Text = "+ Management tools"
$tr3.LanguageID = "sv-SE"
[void]$tr3.InsertBefore(" ")
[void]$tr3.InsertBefore("SQL Server")

$freshTr      = $content3.TextFrame.TextRange
$firstRun3    = $freshTr.Paragraphs(1, 1)
$firstRun3.LanguageID = "en-GB"

$freshTr2 = $content3.TextFrame.TextRange
[void]$freshTr2.InsertBefore("Latest Particular Platform`rMSMQ`r")
